$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "301.83"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-1.00%"
$ws.Range("E2").ClearFormats()
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "2"
$ws.Range("G2").ClearFormats()

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "38.86"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "10.26%"
$ws.Range("E3").ClearFormats()
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "2"
$ws.Range("G3").ClearFormats()

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.991"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-3.45%"
$ws.Range("E4").ClearFormats()
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "2"
$ws.Range("G4").ClearFormats()

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07714"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-1.14%"
$ws.Range("E5").ClearFormats()
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "2"
$ws.Range("G5").ClearFormats()

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.202"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-7.35%"
$ws.Range("E6").ClearFormats()
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "2"
$ws.Range("G6").ClearFormats()

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.002"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.34%"
$ws.Range("E7").ClearFormats()
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "2"
$ws.Range("G7").ClearFormats()

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "4.001"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.41%"
$ws.Range("E8").ClearFormats()
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "2"
$ws.Range("G8").ClearFormats()

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9202"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.67%"
$ws.Range("E9").ClearFormats()
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "2"
$ws.Range("G9").ClearFormats()

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09091"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-8.39%"
$ws.Range("E10").ClearFormats()
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "2"
$ws.Range("G10").ClearFormats()

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1796"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.39%"
$ws.Range("E11").ClearFormats()
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "2"
$ws.Range("G11").ClearFormats()

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08485"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.48%"
$ws.Range("E12").ClearFormats()
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "2"
$ws.Range("G12").ClearFormats()

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03610"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "8.80%"
$ws.Range("E13").ClearFormats()
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "2"
$ws.Range("G13").ClearFormats()

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09956"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.38%"
$ws.Range("E14").ClearFormats()
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "2"
$ws.Range("G14").ClearFormats()

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001470"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.95%"
$ws.Range("E15").ClearFormats()
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "2"
$ws.Range("G15").ClearFormats()

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005744"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.33%"
$ws.Range("E16").ClearFormats()
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "2"
$ws.Range("G16").ClearFormats()

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.483"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.46%"
$ws.Range("E17").ClearFormats()
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "2"
$ws.Range("G17").ClearFormats()

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.216"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.09%"
$ws.Range("E18").ClearFormats()
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "2"
$ws.Range("G18").ClearFormats()

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2.82%"
$ws.Range("E19").ClearFormats()
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "2"
$ws.Range("G19").ClearFormats()

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1319"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "1.85%"
$ws.Range("E20").ClearFormats()
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "2"
$ws.Range("G20").ClearFormats()

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.550"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "4.68%"
$ws.Range("E21").ClearFormats()
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "2"
$ws.Range("G21").ClearFormats()

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-2.69%"
$ws.Range("E22").ClearFormats()
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "2"
$ws.Range("G22").ClearFormats()

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04663"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.00%"
$ws.Range("E23").ClearFormats()
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "2"
$ws.Range("G23").ClearFormats()

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001232"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.48%"
$ws.Range("E24").ClearFormats()
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "2"
$ws.Range("G24").ClearFormats()

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004451"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.32%"
$ws.Range("E25").ClearFormats()
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "2"
$ws.Range("G25").ClearFormats()

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001303"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.29%"
$ws.Range("E26").ClearFormats()
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "2"
$ws.Range("G26").ClearFormats()

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0004762"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "40.26%"
$ws.Range("E27").ClearFormats()
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "2"
$ws.Range("G27").ClearFormats()

$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "2"
$ws.Range("G28").ClearFormats()

$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "2"
$ws.Range("G29").ClearFormats()

$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "2"
$ws.Range("G30").ClearFormats()

$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "2"
$ws.Range("G31").ClearFormats()

$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "2"
$ws.Range("G32").ClearFormats()

$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "2"
$ws.Range("G33").ClearFormats()

$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "2"
$ws.Range("G34").ClearFormats()

$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "2"
$ws.Range("G35").ClearFormats()

$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "2"
$ws.Range("G36").ClearFormats()

$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "2"
$ws.Range("G37").ClearFormats()

$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "2"
$ws.Range("G38").ClearFormats()

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01739"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-4.08%"
$ws.Range("E39").ClearFormats()
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "2"
$ws.Range("G39").ClearFormats()

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04686"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-3.02%"
$ws.Range("E40").ClearFormats()
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "2"
$ws.Range("G40").ClearFormats()

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007831"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.47%"
$ws.Range("E41").ClearFormats()
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "2"
$ws.Range("G41").ClearFormats()

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1390"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-1.81%"
$ws.Range("E42").ClearFormats()
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "2"
$ws.Range("G42").ClearFormats()

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007697"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-14.99%"
$ws.Range("E43").ClearFormats()
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "2"
$ws.Range("G43").ClearFormats()

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002306"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "9.33%"
$ws.Range("E44").ClearFormats()
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "2"
$ws.Range("G44").ClearFormats()

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009472"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "0.49%"
$ws.Range("E45").ClearFormats()
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "2"
$ws.Range("G45").ClearFormats()

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006063"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.72%"
$ws.Range("E46").ClearFormats()
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "2"
$ws.Range("G46").ClearFormats()

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000752"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.27%"
$ws.Range("E47").ClearFormats()
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "2"
$ws.Range("G47").ClearFormats()

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.676"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "101.11%"
$ws.Range("E48").ClearFormats()
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "2"
$ws.Range("G48").ClearFormats()

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "34.96%"
$ws.Range("E49").ClearFormats()
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "2"
$ws.Range("G49").ClearFormats()

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002105"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.27%"
$ws.Range("E50").ClearFormats()
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "2"
$ws.Range("G50").ClearFormats()

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002005"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.27%"
$ws.Range("E51").ClearFormats()
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "2"
$ws.Range("G51").ClearFormats()
